$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "212.50").
# Force Text format so Excel keeps the exact original string (no trailing-zero
# loss / no silent numeric coercion), matching the inline-string cells in the
# source file. The style is reset back to Normal afterwards so no stray
# number-format styling is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.585.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.95%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.580.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("E8").Value = '  +7.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.16'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.26%  '

$ws.Range("E10").Value = '  -1.17%  '

$ws.Range("E11").Value = '  -0.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0883'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.805.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.565.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.562.53'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  -0.48%  '

$ws.Range("E21").Value = '  -1.57%  '

$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.77%  '

$ws.Range("E25").Value = '  +5.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '

$ws.Range("E28").Value = '  -1.59%  '

$ws.Range("E29").Value = '  -1.72%  '

$ws.Range("E30").Value = '  +0.35%  '

$ws.Range("E31").Value = '  -1.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0465'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.75%  '

$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.398.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.28%  '

$ws.Range("E36").Value = '  -2.09%  '

$ws.Range("E37").Value = '  -2.58%  '

$ws.Range("E38").Value = '  +1.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("E41").Value = '  -1.27%  '

$ws.Range("E42").Value = '  +0.39%  '

$ws.Range("E43").Value = '  -1.09%  '

$ws.Range("E44").Value = '  -0.84%  '

$ws.Range("E45").Value = '  +2.50%  '

$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.717.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("E50").Value = '  -2.36%  '
